# Update NATMI ligand-receptor edge table (a -> F11r) with recomputed TPM-based statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.344207
$ws.Range("H2").Value = 1.032621
$ws.Range("I2").Value = 0.6985282229833164
$ws.Range("J2").Value = 0.6985282229833165
$ws.Range("M2").Value = 32.736679
$ws.Range("N2").Value = 98.210037
$ws.Range("O2").Value = 0.8346853755332739
$ws.Range("P2").Value = 0.834685375533274
$ws.Range("Q2").Value = 11.268194068553
$ws.Range("R2").Value = 101.413746616977
$ws.Range("S2").Value = 0.5830512921214199
$ws.Range("T2").Value = 0.5830512921214202

# Row 3: ECs -> FAPs
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.344207
$ws.Range("H3").Value = 1.032621
$ws.Range("I3").Value = 0.6985282229833164
$ws.Range("J3").Value = 0.6985282229833165
$ws.Range("O3").Value = 0.001766029048926899
$ws.Range("P3").Value = 0.0017660290489269
$ws.Range("Q3").Value = 0.02384126838366666
$ws.Range("R3").Value = 0.214571415453
$ws.Range("S3").Value = 0.001233621133283823
$ws.Range("T3").Value = 0.001233621133283824

# Row 4: ECs -> MuSCs
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.344207
$ws.Range("H4").Value = 1.032621
$ws.Range("I4").Value = 0.6985282229833164
$ws.Range("J4").Value = 0.6985282229833165
$ws.Range("M4").Value = 0.5119106666666666
$ws.Range("N4").Value = 1.535732
$ws.Range("O4").Value = 0.01305215923234471
$ws.Range("P4").Value = 0.01305215923234471
$ws.Range("Q4").Value = 0.1762032348413333
$ws.Range("R4").Value = 1.585829113572
$ws.Range("S4").Value = 0.009117301594665039
$ws.Range("T4").Value = 0.00911730159466504

# Row 5: ECs -> Resolving-Mac
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.344207
$ws.Range("H5").Value = 1.032621
$ws.Range("I5").Value = 0.6985282229833164
$ws.Range("J5").Value = 0.6985282229833165
$ws.Range("M5").Value = 5.902527666666667
$ws.Range("N5").Value = 17.707583
$ws.Range("O5").Value = 0.1504964361854544
$ws.Range("P5").Value = 0.1504964361854544
$ws.Range("Q5").Value = 2.031691340560333
$ws.Range("R5").Value = 18.285222065043
$ws.Range("S5").Value = 0.1051260081339476
$ws.Range("T5").Value = 0.1051260081339476

# Row 6: FAPs -> ECs
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1485533333333333
$ws.Range("H6").Value = 0.44566
$ws.Range("I6").Value = 0.3014717770166836
$ws.Range("J6").Value = 0.3014717770166836
$ws.Range("M6").Value = 32.736679
$ws.Range("N6").Value = 98.210037
$ws.Range("O6").Value = 0.8346853755332739
$ws.Range("P6").Value = 0.834685375533274
$ws.Range("Q6").Value = 4.863142787713334
$ws.Range("R6").Value = 43.76828508942
$ws.Range("S6").Value = 0.2516340834118539
$ws.Range("T6").Value = 0.251634083411854

# Row 7: FAPs -> FAPs
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1485533333333333
$ws.Range("H7").Value = 0.44566
$ws.Range("I7").Value = 0.3014717770166836
$ws.Range("J7").Value = 0.3014717770166836
$ws.Range("O7").Value = 0.001766029048926899
$ws.Range("P7").Value = 0.0017660290489269
$ws.Range("Q7").Value = 0.01028944759777778
$ws.Range("R7").Value = 0.09260502838
$ws.Range("S7").Value = 0.0005324079156430759
$ws.Range("T7").Value = 0.000532407915643076

# Row 8: FAPs -> MuSCs
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1485533333333333
$ws.Range("H8").Value = 0.44566
$ws.Range("I8").Value = 0.3014717770166836
$ws.Range("J8").Value = 0.3014717770166836
$ws.Range("M8").Value = 0.5119106666666666
$ws.Range("N8").Value = 1.535732
$ws.Range("O8").Value = 0.01305215923234471
$ws.Range("P8").Value = 0.01305215923234471
$ws.Range("Q8").Value = 0.07604603590222223
$ws.Range("R8").Value = 0.68441432312
$ws.Range("S8").Value = 0.003934857637679673
$ws.Range("T8").Value = 0.003934857637679674

# Row 9: FAPs -> Resolving-Mac
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1485533333333333
$ws.Range("H9").Value = 0.44566
$ws.Range("I9").Value = 0.3014717770166836
$ws.Range("J9").Value = 0.3014717770166836
$ws.Range("M9").Value = 5.902527666666667
$ws.Range("N9").Value = 17.707583
$ws.Range("O9").Value = 0.1504964361854544
$ws.Range("P9").Value = 0.1504964361854544
$ws.Range("Q9").Value = 0.8768401599755556
$ws.Range("R9").Value = 7.89156143978
$ws.Range("S9").Value = 0.04537042805150687
$ws.Range("T9").Value = 0.04537042805150687
